# Weekly update: a new price-report row is added for "Ajo" (Garlic) at
# "Terminal La Palmera de La Serena" on top of the existing historical
# block (rows 146-174), pushing the older rows down by one (146-174 -> 147-175).
#
# The new row duplicates all the static attributes of the former row 146
# (Mercado, Region, Codreg, Categoria, Variedad, Calidad, Unidad, Origen,
# Kg/Unidades, Clasificacion) but carries a new report date (serial 44522)
# while keeping the same price/volume figures as the former row 146.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 146; this shifts the previous rows
# 146:174 down to 147:175 and grows the sheet dimension to A1:R175.
$ws.Rows(146).Insert()

# Populate the newly inserted row 146 with the weekly report values.
$ws.Cells.Item(146, 1).Value = 8
$ws.Cells.Item(146, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(146, 3).Value = "Coquimbo"
$ws.Cells.Item(146, 4).Value = 44522
$ws.Cells.Item(146, 5).Value = 4
$ws.Cells.Item(146, 6).Value = 100112003
$ws.Cells.Item(146, 7).Value = "Ajo"
$ws.Cells.Item(146, 8).Value = "Chino"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 600
$ws.Cells.Item(146, 11).Value = 14000
$ws.Cells.Item(146, 12).Value = 15000
$ws.Cells.Item(146, 13).Value = 14500
$ws.Cells.Item(146, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(146, 15).Value = "China"
$ws.Cells.Item(146, 16).Value = 1450
$ws.Cells.Item(146, 17).Value = 10
$ws.Cells.Item(146, 18).Value = "Hortaliza"

# Match the date cell's number format/style to the rest of column D
# (the insert already carries the style down, but make sure explicitly).
$ws.Cells.Item(146, 4).NumberFormat = $ws.Cells.Item(147, 4).NumberFormat
